$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 131: Mindful Study / Grade 5 Tincture of Mind
$ws.Range("H131").Value = 4831.4287
$ws.Range("I131").Value = 1130
$ws.Range("J131").Value = 4948.316
$ws.Range("K131").Value = 3390
$ws.Range("L131").Value = 14844.948
$ws.Range("M131").Value = 1650
$ws.Range("N131").Value = -24924.948

$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots / Bronze Ingot
$ws.Range("H2").Value = 2014.8572
$ws.Range("I2").Value = 1330.3334
$ws.Range("J2").Value = 2528.25
$ws.Range("K2").Value = 1330.3334
$ws.Range("L2").Value = 2528.25
$ws.Range("M2").Value = -1217.3334
$ws.Range("N2").Value = -2754.25

# Row 17: Cook Intentions / Amateur's Skillet
$ws.Range("H17").Value = 75000
$ws.Range("J17").Value = 75000
$ws.Range("L17").Value = 75000
$ws.Range("N17").Value = -75346

# Row 22: Kiss the Pan (Good-bye) / Initiate's Skillet
$ws.Range("H22").Value = 11966.667
$ws.Range("I22").Value = 1000
$ws.Range("J22").Value = 17450
$ws.Range("K22").Value = 1000
$ws.Range("L22").Value = 17450
$ws.Range("M22").Value = -701
$ws.Range("N22").Value = -18048

# Row 116: No Scope / Titanbronze Ingot
$ws.Range("H116").Value = 2014.8572
$ws.Range("I116").Value = 1330.3334
$ws.Range("J116").Value = 2528.25
$ws.Range("K116").Value = 1330.3334
$ws.Range("L116").Value = 2528.25
$ws.Range("M116").Value = 963.6666
$ws.Range("N116").Value = -7116.25

# Row 135: Forgiveness for My Shins / Ruthenium Sabatons of Fending
$ws.Range("H135").Value = 69857
$ws.Range("J135").Value = 69857
$ws.Range("L135").Value = 69857
$ws.Range("N135").Value = -79997

# Row 139: Backing up My Words / Titanium Gold Thornplate of Fending
$ws.Range("H139").Value = 26560
$ws.Range("J139").Value = 26560
$ws.Range("L139").Value = 26560
$ws.Range("N139").Value = -36840

$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells / Bronze Ingot
$ws.Range("H3").Value = 2014.8572
$ws.Range("I3").Value = 1330.3334
$ws.Range("J3").Value = 2528.25
$ws.Range("K3").Value = 1330.3334
$ws.Range("L3").Value = 2528.25
$ws.Range("M3").Value = -1216.3334
$ws.Range("N3").Value = -2756.25

# Row 75: I Saw the Pine / Hardsilver Saw
$ws.Range("H75").Value = 7360
$ws.Range("I75").Value = 2266.6667
$ws.Range("J75").Value = 15000
$ws.Range("K75").Value = 2266.6667
$ws.Range("L75").Value = 15000
$ws.Range("M75").Value = -1330.6667
$ws.Range("N75").Value = -16872

# Row 78: I Came, I Sawed, I Conquered (L) / Hardsilver Saw
$ws.Range("H78").Value = 7360
$ws.Range("I78").Value = 2266.6667
$ws.Range("J78").Value = 15000
$ws.Range("K78").Value = 6800.000100000001
$ws.Range("L78").Value = 45000
$ws.Range("M78").Value = -2120.000100000001
$ws.Range("N78").Value = -54360

# Row 81: Diamond Sawdust / Titanium Battleaxe
$ws.Range("H81").Value = 21675
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 21675
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 21675
$ws.Range("M81").ClearContents()
$ws.Range("N81").Value = -23797

# Row 84: I'm a Lumberjack and I'm Okay (L) / Titanium Battleaxe
$ws.Range("H84").Value = 21675
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 21675
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 65025
$ws.Range("M84").ClearContents()
$ws.Range("N84").Value = -75633

# Row 135: Axes to the Maxes / Ruthenium War Axe
$ws.Range("H135").Value = 45000
$ws.Range("J135").Value = 45000
$ws.Range("L135").Value = 45000
$ws.Range("N135").Value = -55140

$ws = $wb.Worksheets.Item("CRP")
# Row 22: Driving Up the Wall / Elm Lumber
$ws.Range("H22").Value = 2597.5
$ws.Range("I22").Value = 2597.5
$ws.Range("K22").Value = 2597.5
$ws.Range("M22").Value = -2247.5

# Row 80: The Long Armillae of the Law / Hallowed Chestnut Armillae
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()

# Row 83: Wooden Ambitions (L) / Hallowed Chestnut Armillae
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 5: What a Sap / Maple Syrup
$ws.Range("H5").Value = 1042.4193
$ws.Range("I5").Value = 636.4211
$ws.Range("K5").Value = 1909.2633
$ws.Range("M5").Value = -1797.2633

# Row 31: Food Fight / Shepherd's Pie
$ws.Range("H31").Value = 3000
$ws.Range("J31").Value = 3000
$ws.Range("L31").Value = 9000
$ws.Range("N31").Value = -9576

# Row 46: Feeding Frenzy / Acorn Cookie
$ws.Range("H46").Value = 729.125
$ws.Range("I46").Value = 126.5
$ws.Range("J46").Value = 930
$ws.Range("K46").Value = 379.5
$ws.Range("L46").Value = 2790
$ws.Range("M46").Value = -288.5
$ws.Range("N46").Value = -2972

# Row 54: Good Eats in Ishgard / Salt Cod Puffs
$ws.Range("H54").Value = 4
$ws.Range("I54").Value = 4
$ws.Range("K54").Value = 12
$ws.Range("M54").Value = 547

# Row 57: The Egg Files / Deviled Eggs
$ws.Range("H57").Value = 4600
$ws.Range("I57").Value = 1800
$ws.Range("J57").Value = 13000
$ws.Range("K57").Value = 5400
$ws.Range("L57").Value = 39000
$ws.Range("M57").Value = -4841
$ws.Range("N57").Value = -40118

# Row 58: Bread in the Clouds / La Noscean Toast
$ws.Range("H58").Value = 2876.25
$ws.Range("J58").Value = 3500
$ws.Range("L58").Value = 10500
$ws.Range("N58").Value = -10756

# Row 86: Let's Not Get Sappy / Birch Syrup
$ws.Range("H86").Value = 1003
$ws.Range("J86").Value = 1003
$ws.Range("L86").Value = 3009
$ws.Range("N86").Value = -5381

# Row 89: Luxury Spillover (L) / Birch Syrup
$ws.Range("H89").Value = 1003
$ws.Range("J89").Value = 1003
$ws.Range("L89").Value = 9027
$ws.Range("N89").Value = -20883

# Row 100: Souper / Gameni
$ws.Range("H100").Value = 7185.4287
$ws.Range("J100").Value = 8049.6665
$ws.Range("L100").Value = 24148.9995
$ws.Range("N100").Value = -25770.9995

# Row 114: One Last Meal / Mushroom Saute
$ws.Range("H114").Value = 1387.6
$ws.Range("I114").Value = 229
$ws.Range("J114").Value = 2160
$ws.Range("K114").Value = 687
$ws.Range("L114").Value = 6480
$ws.Range("M114").Value = 2567
$ws.Range("N114").Value = -12988

# Row 122: Salt of the North / Northern Sea Salt
$ws.Range("H122").Value = 675.2105
$ws.Range("J122").Value = 865.36365
$ws.Range("L122").Value = 7788.27285
$ws.Range("N122").Value = -12688.27285

# Row 135: Not-so-secret Ingredient / Royal Maple Syrup
$ws.Range("H135").Value = 1042.4193
$ws.Range("I135").Value = 636.4211
$ws.Range("K135").Value = 5727.7899
$ws.Range("M135").Value = -3192.7899

$ws = $wb.Worksheets.Item("GSM")
# Row 24: Bad Guys Eat Brass / Brass Ring of Crafting
$ws.Range("H24").Value = 118528.14
$ws.Range("J24").Value = 37056.285
$ws.Range("L24").Value = 37056.285
$ws.Range("N24").Value = -37402.285

$ws = $wb.Worksheets.Item("LTW")
# Row 133: The Perfect Accessory / Loboskin Amulet of Fending
$ws.Range("H133").Value = 22020.25
$ws.Range("J133").Value = 22020.25
$ws.Range("L133").Value = 22020.25
$ws.Range("N133").Value = -27080.25

$ws = $wb.Worksheets.Item("WVR")
# Row 15: Workplace Safety / Cotton Scarf
$ws.Range("H15").Value = 70007
$ws.Range("J15").Value = 70007
$ws.Range("L15").Value = 70007
$ws.Range("N15").Value = -70583

# Row 18: Welcome to the Cotton Club / Cotton Halfgloves
$ws.Range("H18").Value = 36669
$ws.Range("J18").Value = 36669
$ws.Range("L18").Value = 36669
$ws.Range("N18").Value = -37015

# Row 20: Read the Fine Print / Cotton Shepherd's Tunic
$ws.Range("H20").Value = 48340.668
$ws.Range("J20").Value = 48340.668
$ws.Range("L20").Value = 48340.668
$ws.Range("N20").Value = -48820.668

# Row 22: Better Shroud than Sorry / Cotton Kurta
$ws.Range("H22").Value = 34000
$ws.Range("J22").Value = 34000
$ws.Range("L22").Value = 34000
$ws.Range("N22").Value = -34586

# Row 31: Whatchoo Talking About / Cotton Doublet Vest of Crafting
$ws.Range("H31").Value = 70019
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()

# Row 69: Fashion Patrol / Holy Rainbow Sarouel of Casting
$ws.Range("H69").Value = 25000
$ws.Range("J69").Value = 25000
$ws.Range("L69").Value = 25000
$ws.Range("N69").Value = -26498

# Row 72: Dress Code Violation (L) / Holy Rainbow Sarouel of Casting
$ws.Range("H72").Value = 25000
$ws.Range("J72").Value = 25000
$ws.Range("L72").Value = 75000
$ws.Range("N72").Value = -82488

# Row 75: Storm upon Bald Mountain / Ramie Turban of Crafting
$ws.Range("H75").Value = 26559
$ws.Range("I75").Value = 8118
$ws.Range("J75").Value = 45000
$ws.Range("K75").Value = 8118
$ws.Range("L75").Value = 45000
$ws.Range("M75").Value = -7182
$ws.Range("N75").Value = -46872

# Row 78: Abrupt Apprentices (L) / Ramie Turban of Crafting
$ws.Range("H78").Value = 26559
$ws.Range("I78").Value = 8118
$ws.Range("J78").Value = 45000
$ws.Range("K78").Value = 24354
$ws.Range("L78").Value = 135000
$ws.Range("M78").Value = -19674
$ws.Range("N78").Value = -144360

# Row 101: Who War It Better / Serge Hose of Aiming
$ws.Range("H101").Value = 39251
$ws.Range("J101").Value = 39251
$ws.Range("L101").Value = 39251
$ws.Range("N101").Value = -45741
